# Apply updated vm_pu values for the 380 kV case (rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colLetters = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
    2 = @(1.02,1.06138285887833,1.074568748949277,1.056788327412819,1.081601295330561,1.059842832034928,1.066358825236371,1.077257616042017,1.059525077788409,1.084271692871229,1.067873177081273)
    3 = @(1.02,1.062874519181515,1.075585571559972,1.058100983930602,1.082811667586894,1.06036644463045,1.067501920918978,1.078091098302245,1.060650173606004,1.085299591951552,1.069017896090865)
    4 = @(1.02,1.063837731300729,1.076242451908338,1.058948094291611,1.083593901699216,1.060703197400707,1.068239163609764,1.078628733634291,1.061375378027037,1.085963179335528,1.069756185750835)
    5 = @(1.02,1.064242197449115,1.076518351535174,1.059303684464307,1.083922528402925,1.060844278183804,1.068548528549638,1.078854356195668,1.061679588731147,1.086241789202356,1.070065990024408)
    6 = @(1.02,1.064310081887757,1.076564661576033,1.059363358462173,1.083977693244128,1.060867937615897,1.068600438941396,1.078892215905949,1.061730628238707,1.086288547911638,1.070117974134877)
    7 = @(1.02,1.063843137626342,1.07624613948166,1.058952847797893,1.083598293701644,1.060705084451038,1.068243299598128,1.078631749978562,1.061379445510669,1.085966903552114,1.069760327612776)
    8 = @(1.02,1.061887390092151,1.074912613108753,1.057232419587865,1.082010547464824,1.060020217620864,1.066745644182135,1.077539646063661,1.059905894764549,1.084619396123354,1.068260545354309)
    9 = @(1.02,1.058425470108931,1.072554391227226,1.054183111412894,1.079205186842102,1.058797486029985,1.064087775430702,1.075602160263455,1.057287492543463,1.082232978928671,1.065598902124387)
    10 = @(1.02,1.056106461456617,1.07097637145905,1.05213786852602,1.077329568389204,1.057971445480878,1.062302795502224,1.074301488019853,1.055526782228684,1.080633721543376,1.063811387319671)
    11 = @(1.02,1.055099559513634,1.070291624697776,1.051249215338701,1.076516061067768,1.057611135964775,1.061526685571158,1.073736091338025,1.054760693070262,1.07993918916287,1.063034175223482)
    12 = @(1.02,1.054725126676354,1.070037056091662,1.05091866257615,1.07621367910519,1.05747690230624,1.061237914513326,1.073525743102794,1.054475570487988,1.079680896253503,1.062744994077625)
    13 = @(1.02,1.054805463176939,1.070091672097645,1.05098958856599,1.076278550622401,1.057505713982635,1.061299879159064,1.073570878735929,1.054536755863384,1.079736315199371,1.062807046720267)
    14 = @(1.02,1.055068617502212,1.070270586557677,1.05122190133284,1.076491070381212,1.057600048323119,1.061502825694332,1.073718710742837,1.054737136282752,1.079917845001518,1.063010281462894)
    15 = @(1.02,1.05523069898242,1.070380791946054,1.051364974682636,1.076621982928201,1.057658117884198,1.061627802671866,1.07380975042588,1.054860522417061,1.080029649903106,1.063135435922071)
    16 = @(1.02,1.056173227257749,1.07102178481201,1.052196780474085,1.077383529193789,1.057995302337199,1.062354235202074,1.074338964850503,1.055577546568701,1.080679771836157,1.063862900069795)
    17 = @(1.02,1.056763705069859,1.071423470139454,1.052717727324809,1.077860860834965,1.058206102867714,1.062809044121921,1.074670335600207,1.056026322928458,1.081087025095611,1.064318354870476)
    18 = @(1.02,1.057107855603511,1.071657626449996,1.053021293203536,1.078139150249588,1.058328805893823,1.063074018098843,1.07486340688932,1.056287731141259,1.081324372261266,1.064583705140838)
    19 = @(1.02,1.057225157380613,1.071737444106474,1.053124751801058,1.078234017868217,1.058370601594096,1.063164315320104,1.074929203447747,1.056376804555131,1.081405268305769,1.06467413059451)
    20 = @(1.02,1.056700379894679,1.071380387614987,1.052661865093865,1.077809661162908,1.058183512198905,1.062760279384986,1.074634804581861,1.055978210299728,1.081043351069023,1.064269520882021)
    21 = @(1.02,1.054991136919473,1.07021790690756,1.051153503990455,1.076428494430228,1.057572280246426,1.061443076570806,1.073675187179316,1.054678144853294,1.079864397675981,1.062950447488759)
    22 = @(1.02,1.053914005975993,1.069485715560506,1.050202428621087,1.075558887169891,1.057185666039645,1.06061206396934,1.073069898628159,1.053857480699328,1.079121331090725,1.062118254754081)
    23 = @(1.02,1.054485250397692,1.069873988264085,1.0507068715459,1.076019999412129,1.057390837620217,1.061052870955914,1.073390958902527,1.054292842321882,1.079515418388637,1.062559687737137)
    24 = @(1.02,1.056728994635082,1.071399855186018,1.052687107738839,1.077832796492993,1.058193720732697,1.062782315023536,1.074650860173859,1.055999951424086,1.081063086089946,1.064291587813705)
    25 = @(1.02,1.059322368977686,1.073165065396333,1.054973576090396,1.079931365366812,1.059115496179913,1.064777170316077,1.076104618395374,1.057967042849954,1.082851368933092,1.066289276029575)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $colLetters.Length; $i++) {
        $ws.Range("$($colLetters[$i])$r").Value = $vals[$i]
    }
}

Write-Host "Updated $($data.Keys.Count) rows"
